$d = $word.ActiveDocument

# ---- Locate insertion point right before "<lb/>" in "...foeu<del>t</del><lb/>..." ----
$search = $d.Content
$found = $search.Find.Execute("foeu<del>t</del><lb/>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$scoped = $d.Range($search.Start, $search.End)
$scoped.Find.Execute("<lb/>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insertAt = $scoped.Start

$openTag = "<comment>"
$id = "c_68r_03"
$closeTag = "</comment>"

# ---- Insert combined text, then restyle the open/close tag runs (Courier New, blue) ----
$ip = $d.Range($insertAt, $insertAt)
$ip.InsertBefore($openTag + $id + $closeTag)

$openRange = $d.Range($insertAt, $insertAt + $openTag.Length)
$openRange.Font.Name = "Courier New"
$openRange.Font.Color = 16711680   # RGB(0,0,255) -> w:val="0000ff"
$openRange.Font.Size = 9           # -> w:sz/w:szCs val="18" (half-points)

$closeRange = $d.Range($insertAt + $openTag.Length + $id.Length, $insertAt + $openTag.Length + $id.Length + $closeTag.Length)
$closeRange.Font.Name = "Courier New"
$closeRange.Font.Color = 16711680
$closeRange.Font.Size = 9

# ---- Build the middle "c_68r_03" run (no rFonts override, dark-red, smaller) ----
# Stage it at a spot with no ambient rFonts override so the new run doesn't
# inherit Courier New, then cut/paste it into the real target position.
$stage = $d.Content
$stage.Find.Execute("violent ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$stageStart = $stage.Start
$stageIp = $d.Range($stageStart, $stageStart)
$stageIp.InsertBefore($id)
$stageRange = $d.Range($stageStart, $stageStart + $id.Length)
$stageRange.Font.Color = 1118633   # RGB(169,17,17) -> w:val="a91111"
$stageRange.Font.Size = 8          # -> w:sz val="16" (half-points)
$stageRange.Cut()

$idTarget = $d.Range($insertAt + $openTag.Length, $insertAt + $openTag.Length + $id.Length)
$idTarget.Paste()

# ---- sectPr: add footer distance (w:footer="720" twips = 36pt) ----
$d.PageSetup.FooterDistance = 36
